$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A=465; B="Monday, Jan 16"; C="2:09 AM";  D="QY5102"; E="Leipzig";    F="(LEJ)"; G="DHL ";       H="B738"; I="(EC-IXO)"; J="1:53 AM";  L="0 hours, -16 minutes" },
    @{ A=466; B="Monday, Jan 16"; C="6:04 AM";  D="UNKNOWN"; E="Katowice";   F="(KTW)"; G="Enter Air "; H="B738"; I="(SP-ESD)"; J="6:04 AM";  L="0 hours, 0 minutes" },
    @{ A=467; B="Monday, Jan 16"; C="7:30 AM";  D="FR4105"; E="Wroclaw";    F="(WRO)"; G="Ryanair ";   H="B738"; I="(SP-RKI)"; J="7:29 AM";  L="0 hours, -1 minutes" },
    @{ A=468; B="Monday, Jan 16"; C="7:52 AM";  D="P81956"; E="Cologne";    F="(CGN)"; G="SprintAir "; H="SF34"; I="(SP-KPE)"; J="7:34 AM";  L="0 hours, -18 minutes" },
    @{ A=469; B="Monday, Jan 16"; C="9:25 AM";  D="FR6098"; E="Gothenburg"; F="(GOT)"; G="Ryanair ";   H="B738"; I="(SP-RSO)"; J="9:33 AM";  L="0 hours, 8 minutes" },
    @{ A=470; B="Monday, Jan 16"; C="9:25 AM";  D="FR6845"; E="Copenhagen"; F="(CPH)"; G="Ryanair ";   H="B738"; I="(SP-RSL)"; J="9:22 AM";  L="0 hours, -3 minutes" },
    @{ A=471; B="Monday, Jan 16"; C="9:30 AM";  D="FR6112"; E="Lublin";     F="(LUZ)"; G="Ryanair ";   H="B738"; I="(SP-RSW)"; J="9:16 AM";  L="0 hours, -14 minutes" },
    @{ A=472; B="Monday, Jan 16"; C="9:45 AM";  D="FR3278"; E="Oslo";       F="(TRF)"; G="Ryanair ";   H="B738"; I="(SP-RKM)"; J="9:53 AM";  L="0 hours, 8 minutes" },
    @{ A=473; B="Monday, Jan 16"; C="10:35 AM"; D="W61642"; E="Eindhoven";  F="(EIN)"; G="Wizz Air ";  H="A320"; I="(HA-LYH)"; J="10:11 AM"; L="0 hours, -24 minutes" }
)

$startRow = 466
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 12).Value = $row.L
}
